$wb = $excel.ActiveWorkbook

# Add the new "DESI" entry to the ms_source list sheet, inserted before the
# existing "nanoDESI" entry (so it becomes row 7, pushing nanoDESI to row 8)
$msSourceWs = $wb.Worksheets.Item("ms_source list")
$msSourceWs.Range("A7").EntireRow.Insert()
$msSourceWs.Range("A7").Value = "DESI"

# Update the data validation on the "Export as TSV" sheet (column Q) so the
# allowed list range grows from $A$1:$A$7 to $A$1:$A$8
$tsvWs = $wb.Worksheets.Item("Export as TSV")
$qRange = $tsvWs.Range("Q2:Q1048576")
$qRange.Validation.Modify(3, 1, 1, "='ms_source list'!`$A`$1:`$A`$8")
